$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Exported:" timestamp string (row 4, col A)
$ws.Range("A4").Value = "Exported: 2018-03-31 18:54:47"

# 2. Insert two new rows below the current "Logout" row (row 12), pushing the
#    blank spacer row (old row 15) down to row 17.
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(14).Insert()

# 3. Copy the formatting of row 12 (A:U) into the two newly inserted rows so
#    they pick up the same cell styles used by the other data rows.
$ws.Range("A12:U12").Copy($ws.Range("A13:U13"))
$ws.Range("A12:U12").Copy($ws.Range("A14:U14"))

# 4. Populate the three rows with their final content.
#    Row 12 becomes "My Page", row 13 becomes "Admin", row 14 becomes the
#    re-inserted "Logout" entry.
$ws.Range("C12").Value = "My Page"
$ws.Range("K12").Value = "/mypage/"
$ws.Range("R12").Value = 1

$ws.Range("C13").Value = "Admin"
$ws.Range("K13").Value = "/admin/"
$ws.Range("R13").Value = 1

$ws.Range("C14").Value = "Logout"
$ws.Range("K14").Value = "/logout.html"
$ws.Range("R14").Value = 0
